# Generate Report for Handoff
# Updates the localization-status workbook: refreshes the in-flight
# handoff id (24ef0adf-... -> 1a420ffc-...), its handoff timestamps, and
# appends two new rows (for the two .png assets that are now part of the
# handoff) to the Overview sheet and each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
#   Columns: A=File Name  B=zh-cn  C=de-de  D=Latest Handoff Date
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Clear & rebuild the hyperlinks collection so ids/rels are regenerated
# cleanly (per-item Delete/replace does not update the underlying part).
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "1a420ffc-7661-417f-a439-a5077600dcd7.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-24 17:17:42"

$ov.Range("A3").Value = "adfa4149-a733-4de0-911a-9119bcecf0d7.png"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-24 17:17:42"
$ov.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Range("A4").Value = "f8a2f5de-4fe7-476c-856f-96cfc1463034.png"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-24 17:17:42"
$ov.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/1a420ffc-7661-417f-a439-a5077600dcd7.md", "", "", "1a420ffc-7661-417f-a439-a5077600dcd7.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/adfa4149-a733-4de0-911a-9119bcecf0d7.png", "", "", "adfa4149-a733-4de0-911a-9119bcecf0d7.png")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/f8a2f5de-4fe7-476c-856f-96cfc1463034.png", "", "", "f8a2f5de-4fe7-476c-856f-96cfc1463034.png")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
#   Columns: A=Source File Name  B=File Extension  C=Status
#            D=Latest Handoff File  E=Latest Handoff Datetime
#            F=Latest Target File  G=Latest Handback File
#            H=Latest Handback DateTime  I=Reference Tokens
#            J=Handoff Reason  K=Dependency From  L=Error Detail
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = "1a420ffc-7661-417f-a439-a5077600dcd7.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-24 17:17:38"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "adfa4149-a733-4de0-911a-9119bcecf0d7.png"
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "223e671941029660bb1645d0d6490f1bfd0341a4.png"
$zh.Range("E3").Value = "2016-03-24 17:17:38"
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("J3").Value = "IsDependency"
$zh.Range("K3").Value = "e2e\1a420ffc-7661-417f-a439-a5077600dcd7.md"

$zh.Range("A4").Value = "f8a2f5de-4fe7-476c-856f-96cfc1463034.png"
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "cbbbbe6b2189f7571d6789cae734c0824b6f797e.png"
$zh.Range("E4").Value = "2016-03-24 17:17:38"
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("J4").Value = "IsDependency"
$zh.Range("K4").Value = "e2e\1a420ffc-7661-417f-a439-a5077600dcd7.md"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/1a420ffc-7661-417f-a439-a5077600dcd7.md", "", "", "1a420ffc-7661-417f-a439-a5077600dcd7.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.zh-cn.xlf", "", "", "1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/adfa4149-a733-4de0-911a-9119bcecf0d7.png", "", "", "adfa4149-a733-4de0-911a-9119bcecf0d7.png")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/223e671941029660bb1645d0d6490f1bfd0341a4.png", "", "", "223e671941029660bb1645d0d6490f1bfd0341a4.png")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/f8a2f5de-4fe7-476c-856f-96cfc1463034.png", "", "", "f8a2f5de-4fe7-476c-856f-96cfc1463034.png")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/cbbbbe6b2189f7571d6789cae734c0824b6f797e.png", "", "", "cbbbbe6b2189f7571d6789cae734c0824b6f797e.png")

# ---------------------------------------------------------------------
# Sheet 3: "de-de" (same column layout as "zh-cn")
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = "1a420ffc-7661-417f-a439-a5077600dcd7.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.de-de.xlf"
$de.Range("E2").Value = "2016-03-24 17:17:42"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "adfa4149-a733-4de0-911a-9119bcecf0d7.png"
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "223e671941029660bb1645d0d6490f1bfd0341a4.png"
$de.Range("E3").Value = "2016-03-24 17:17:42"
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("J3").Value = "IsDependency"
$de.Range("K3").Value = "e2e\1a420ffc-7661-417f-a439-a5077600dcd7.md"

$de.Range("A4").Value = "f8a2f5de-4fe7-476c-856f-96cfc1463034.png"
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "cbbbbe6b2189f7571d6789cae734c0824b6f797e.png"
$de.Range("E4").Value = "2016-03-24 17:17:42"
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("J4").Value = "IsDependency"
$de.Range("K4").Value = "e2e\1a420ffc-7661-417f-a439-a5077600dcd7.md"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/1a420ffc-7661-417f-a439-a5077600dcd7.md", "", "", "1a420ffc-7661-417f-a439-a5077600dcd7.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.de-de.xlf", "", "", "1a420ffc-7661-417f-a439-a5077600dcd7.a8cf557d24e897810ddb51255d0dddeb8d7834cf.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/adfa4149-a733-4de0-911a-9119bcecf0d7.png", "", "", "adfa4149-a733-4de0-911a-9119bcecf0d7.png")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/223e671941029660bb1645d0d6490f1bfd0341a4.png", "", "", "223e671941029660bb1645d0d6490f1bfd0341a4.png")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/f8a2f5de-4fe7-476c-856f-96cfc1463034.png", "", "", "f8a2f5de-4fe7-476c-856f-96cfc1463034.png")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/cbbbbe6b2189f7571d6789cae734c0824b6f797e.png", "", "", "cbbbbe6b2189f7571d6789cae734c0824b6f797e.png")

Write-Output "Report generated for handoff"
